$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = 'summ47305199'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'6264.443779274492'
$arr[0,2] = [double]'0.3117617096370295'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'174.1211657367978'
$arr[1,2] = [double]'0.9442604047796873'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-83.3680483818855'
$arr[2,2] = [double]'0.9644044904744157'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'236.7816179125396'
$arr[3,2] = [double]'0.7803179106390103'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'697.590589402412'
$arr[4,2] = [double]'0.3862832285418086'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-230.8192663111598'
$arr[5,2] = [double]'0.7874424704744367'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1477.635888804815'
$arr[6,2] = [double]'0.06704358511993058'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'-103.1638084882058'
$arr[7,2] = [double]'0.7020090953559899'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1436.020146987179'
$arr[8,2] = [double]'0.01613545524123585'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-6.435834201163505'
$arr[9,2] = [double]'0.8217432300920015'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'464.1592549163218'
$arr[10,2] = [double]'0.1217413496397614'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'540.9357211828204'
$arr[11,2] = [double]'0.004830371414727553'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.3250516826446449'
$arr[12,2] = [double]'0.06504753503677504'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'8.656654442691398e-06'
$arr[13,2] = [double]'0.9761264556850359'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'14.33252490871436'
$arr[14,2] = [double]'0.7669818712023923'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-8.731508033447295'
$arr[15,2] = [double]'0.8319781954046517'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'2520.976286310312'
$arr[16,2] = [double]'0.5912699190298034'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-3962.213950875906'
$arr[17,2] = [double]'0.3946669647775711'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-7124.426277433507'
$arr[18,2] = [double]'0.1509513510259737'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = 'summ47476983'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'9458.748584398361'
$arr[0,2] = [double]'0.1138466152542727'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2242.1033375741'
$arr[1,2] = [double]'0.2820280078817987'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-615.1012908341813'
$arr[2,2] = [double]'0.7306173183615793'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'91.20642531794687'
$arr[3,2] = [double]'0.9125305891017931'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'919.9611715347119'
$arr[4,2] = [double]'0.2350449495113479'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'57.2391033044039'
$arr[5,2] = [double]'0.9449306372688479'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1348.972053123113'
$arr[6,2] = [double]'0.0881593409933355'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'139.0056765908662'
$arr[7,2] = [double]'0.6042320423726015'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1541.918743193851'
$arr[8,2] = [double]'0.007551481989724755'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-15.95358858213493'
$arr[9,2] = [double]'0.5634517385941555'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'301.171880177191'
$arr[10,2] = [double]'0.3107652687577218'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'498.4942756054291'
$arr[11,2] = [double]'0.007844752695359246'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.2856325698350092'
$arr[12,2] = [double]'0.09164733144404129'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'5.395577034850499e-05'
$arr[13,2] = [double]'0.8514123580227102'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'4.380723201062377'
$arr[14,2] = [double]'0.9274712894431764'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-25.47028513667758'
$arr[15,2] = [double]'0.5228119087516211'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'1463.812631483728'
$arr[16,2] = [double]'0.7554680847595993'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-3954.172611528266'
$arr[17,2] = [double]'0.381218759184577'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-7849.056281420352'
$arr[18,2] = [double]'0.1002759343696318'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = 'summ47654368'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'3526.66909159663'
$arr[0,2] = [double]'0.5519836912429669'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2571.438565985788'
$arr[1,2] = [double]'0.249812515715047'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-765.0892539285911'
$arr[2,2] = [double]'0.6457565941654559'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'-201.958168635979'
$arr[3,2] = [double]'0.8032810172000017'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'773.6415571989753'
$arr[4,2] = [double]'0.3270947744611918'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-109.2000492935445'
$arr[5,2] = [double]'0.8957245090109434'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1271.502477953967'
$arr[6,2] = [double]'0.1117415678092454'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'-17.35159060370108'
$arr[7,2] = [double]'0.9469246476782768'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1188.657919349037'
$arr[8,2] = [double]'0.03876511994689568'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'8.846344225698942'
$arr[9,2] = [double]'0.7459911573261357'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'236.916425058678'
$arr[10,2] = [double]'0.429870603993504'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'595.0297275595816'
$arr[11,2] = [double]'0.001570193590330965'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.1673360429507389'
$arr[12,2] = [double]'0.3274176802369223'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-8.688150551975187e-05'
$arr[13,2] = [double]'0.7621278089431921'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'52.45570038185288'
$arr[14,2] = [double]'0.279197699987609'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'14.14424788826765'
$arr[15,2] = [double]'0.7220429777014392'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'1232.366883319334'
$arr[16,2] = [double]'0.7900181518447658'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-6046.065548233317'
$arr[17,2] = [double]'0.1768740614893659'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-6892.615562659079'
$arr[18,2] = [double]'0.1549626925273081'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = 'summ47829878'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'10999.28664934736'
$arr[0,2] = [double]'0.06141151159960664'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2661.060526673187'
$arr[1,2] = [double]'0.2192804342040439'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-2386.565172651282'
$arr[2,2] = [double]'0.1729678177231806'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'-186.8170801311937'
$arr[3,2] = [double]'0.8204945290203544'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'684.1507505290646'
$arr[4,2] = [double]'0.3790927328861378'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-493.6398292681063'
$arr[5,2] = [double]'0.5487724722051746'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'221.6770745385281'
$arr[6,2] = [double]'0.7771158114124385'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'138.7592645882799'
$arr[7,2] = [double]'0.5893212825646368'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-773.0156301441118'
$arr[8,2] = [double]'0.1744362151690799'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-20.07198080204572'
$arr[9,2] = [double]'0.4635777163562066'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'392.2748555834979'
$arr[10,2] = [double]'0.188175162392914'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'394.4879301699053'
$arr[11,2] = [double]'0.0306469353404468'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.2716554375747139'
$arr[12,2] = [double]'0.1048087873660971'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-1.010208985035391e-05'
$arr[13,2] = [double]'0.9703150653888037'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'10.6989918184643'
$arr[14,2] = [double]'0.8209350309064598'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-29.33285117073705'
$arr[15,2] = [double]'0.460941104283437'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'822.9379604390688'
$arr[16,2] = [double]'0.8545222214576721'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-5860.249654608447'
$arr[17,2] = [double]'0.190517291913683'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-7921.357039189117'
$arr[18,2] = [double]'0.09992526255616298'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = 'summ48009910'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'9324.053023873857'
$arr[0,2] = [double]'0.1160933218316672'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2886.219017818962'
$arr[1,2] = [double]'0.196319135820986'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-629.270027869116'
$arr[2,2] = [double]'0.7115694076597123'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'-489.2859892938632'
$arr[3,2] = [double]'0.5447497984579013'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'33.53391447877021'
$arr[4,2] = [double]'0.9659079229399252'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-775.0030912084367'
$arr[5,2] = [double]'0.3554067525253517'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'695.703273085074'
$arr[6,2] = [double]'0.376893054883124'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'141.9544699255349'
$arr[7,2] = [double]'0.5805520017501525'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1563.319619392382'
$arr[8,2] = [double]'0.006497189063604992'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-9.828160266194445'
$arr[9,2] = [double]'0.7228742442964564'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'511.9485020243777'
$arr[10,2] = [double]'0.07953117664067243'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'332.4026103473077'
$arr[11,2] = [double]'0.08039191975684763'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.1981456817648541'
$arr[12,2] = [double]'0.2393735980414387'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-6.716304922632929e-05'
$arr[13,2] = [double]'0.8057790330935291'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'7.839147013386029'
$arr[14,2] = [double]'0.8678997470273342'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-16.29017507677502'
$arr[15,2] = [double]'0.6862604664083817'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'2887.195008421633'
$arr[16,2] = [double]'0.5275809595499744'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-1045.664403176253'
$arr[17,2] = [double]'0.8157427021552379'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-9337.891829243294'
$arr[18,2] = [double]'0.05726051167211738'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = 'summ48192640'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'7289.723500844177'
$arr[0,2] = [double]'0.2306942602702995'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2516.817808172233'
$arr[1,2] = [double]'0.2623289790193089'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-725.2789024486631'
$arr[2,2] = [double]'0.6896668388499785'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'-194.5146485888953'
$arr[3,2] = [double]'0.8233870250841842'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'734.248999487206'
$arr[4,2] = [double]'0.3657868430133949'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-365.9283087541921'
$arr[5,2] = [double]'0.6733521139737615'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1236.75416087398'
$arr[6,2] = [double]'0.1328234468883536'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'45.87575908013807'
$arr[7,2] = [double]'0.8639972647186006'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-900.0032323778331'
$arr[8,2] = [double]'0.1329572293600494'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-26.7854399081226'
$arr[9,2] = [double]'0.3427590640366711'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'586.5995510178111'
$arr[10,2] = [double]'0.05913367792768713'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'516.4079879595702'
$arr[11,2] = [double]'0.008116801201066868'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.2349856785890758'
$arr[12,2] = [double]'0.1723225231661712'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-7.741722681953796e-05'
$arr[13,2] = [double]'0.784488199546448'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'27.34471858028423'
$arr[14,2] = [double]'0.5827936835217049'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-6.905011518683018'
$arr[15,2] = [double]'0.8639435670017555'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'1106.099028260238'
$arr[16,2] = [double]'0.8142364039319251'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-6025.055008903573'
$arr[17,2] = [double]'0.2001658348216329'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-5062.265559249135'
$arr[18,2] = [double]'0.3082895857803247'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = 'summ48378127'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'6684.49704144878'
$arr[0,2] = [double]'0.2772443993513913'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-1741.656304453086'
$arr[1,2] = [double]'0.4132045535369808'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-353.58581035668'
$arr[2,2] = [double]'0.8451066639403833'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'126.1132372534765'
$arr[3,2] = [double]'0.8802526756827227'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'844.1602346552584'
$arr[4,2] = [double]'0.2990514216473436'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-157.7371600958612'
$arr[5,2] = [double]'0.854529304275494'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'846.2487277751791'
$arr[6,2] = [double]'0.2972474140118442'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'52.75310659349998'
$arr[7,2] = [double]'0.8437772571156869'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1156.981066355245'
$arr[8,2] = [double]'0.05086189049386818'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-13.27719424078026'
$arr[9,2] = [double]'0.6404110107650208'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'276.7528968519009'
$arr[10,2] = [double]'0.3758487496419283'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'475.449860972137'
$arr[11,2] = [double]'0.01369288389098398'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.2167180527959116'
$arr[12,2] = [double]'0.2148663122584808'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-6.088412667224479e-05'
$arr[13,2] = [double]'0.8371707906895633'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'9.555858273727988'
$arr[14,2] = [double]'0.8467579340194774'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'5.371895490326711'
$arr[15,2] = [double]'0.8981855465690306'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'673.2937028061542'
$arr[16,2] = [double]'0.8914985585396057'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-2764.2927473051'
$arr[17,2] = [double]'0.5519559272373014'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-7806.325021587754'
$arr[18,2] = [double]'0.1168552326389852'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = 'summ48559373'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'2766.206570967248'
$arr[0,2] = [double]'0.6413676980215117'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-2252.088136402665'
$arr[1,2] = [double]'0.2776949610930483'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-2001.510508037287'
$arr[2,2] = [double]'0.2835263219242073'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'-749.2803631604525'
$arr[3,2] = [double]'0.3681724394659536'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'523.6769562518116'
$arr[4,2] = [double]'0.5031574136298858'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'-757.3466334381808'
$arr[5,2] = [double]'0.3624165981467464'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1030.222772934864'
$arr[6,2] = [double]'0.1901203808731933'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'181.7486738610716'
$arr[7,2] = [double]'0.4899766405266476'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-818.4016611751243'
$arr[8,2] = [double]'0.1552975409327697'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-3.634991321778152'
$arr[9,2] = [double]'0.8938182590653402'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'347.4686620901791'
$arr[10,2] = [double]'0.2487690877966206'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'562.4349239499097'
$arr[11,2] = [double]'0.00250915734521565'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.1994267683316536'
$arr[12,2] = [double]'0.2421513096587156'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'0.0001999813970620377'
$arr[13,2] = [double]'0.5358504963217419'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'18.73294904687641'
$arr[14,2] = [double]'0.693891839276508'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'25.07673528276289'
$arr[15,2] = [double]'0.5397000667312704'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'-874.4695143096681'
$arr[16,2] = [double]'0.8615405449223907'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-3612.787583713088'
$arr[17,2] = [double]'0.4155526716880218'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-4851.604412804874'
$arr[18,2] = [double]'0.3148492670126323'
$ws.Range("A2:C20").Value = $arr

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = 'summ48738292'
$ws.Range("A2:C17").Clear()
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = 'Intercept'
$arr[0,1] = [double]'8939.254863099552'
$arr[0,2] = [double]'0.1290355802833208'
$arr[1,0] = 'Education[T.Primary/None]'
$arr[1,1] = [double]'-1815.482893801636'
$arr[1,2] = [double]'0.4014389137583351'
$arr[2,0] = 'Education[T.Secondary]'
$arr[2,1] = [double]'-810.9007893392736'
$arr[2,2] = [double]'0.6346379410589764'
$arr[3,0] = 'Education[T.University]'
$arr[3,1] = [double]'41.65623485639006'
$arr[3,2] = [double]'0.959716720326196'
$arr[4,0] = 'Season[T.Spring]'
$arr[4,1] = [double]'798.8048313378027'
$arr[4,2] = [double]'0.3027293011418784'
$arr[5,0] = 'Season[T.Summer]'
$arr[5,1] = [double]'103.3701102577679'
$arr[5,2] = [double]'0.9001235428906171'
$arr[6,0] = 'Season[T.Winter]'
$arr[6,1] = [double]'1222.726265135708'
$arr[6,2] = [double]'0.1216701163138527'
$arr[7,0] = 'HHSize'
$arr[7,1] = [double]'90.10635853356118'
$arr[7,2] = [double]'0.7263907378613678'
$arr[8,0] = 'Sex'
$arr[8,1] = [double]'-1168.245702595202'
$arr[8,2] = [double]'0.0427677218492331'
$arr[9,0] = 'Age'
$arr[9,1] = [double]'-23.42916743451748'
$arr[9,2] = [double]'0.3881483729959176'
$arr[10,0] = 'DistSubcenter_res'
$arr[10,1] = [double]'417.4436520925403'
$arr[10,2] = [double]'0.1522769971610238'
$arr[11,0] = 'DistCenter_res'
$arr[11,1] = [double]'604.3912256868311'
$arr[11,2] = [double]'0.001414369777523669'
$arr[12,0] = 'UrbPopDensity_res'
$arr[12,1] = [double]'0.2381409909862235'
$arr[12,2] = [double]'0.151674976224007'
$arr[13,0] = 'UrbBuildDensity_res'
$arr[13,1] = [double]'-0.0002825068068225374'
$arr[13,2] = [double]'0.3095162295488683'
$arr[14,0] = 'IntersecDensity_res'
$arr[14,1] = [double]'7.612410067778923'
$arr[14,2] = [double]'0.8699213152367583'
$arr[15,0] = 'street_length_res'
$arr[15,1] = [double]'-20.66269174960851'
$arr[15,2] = [double]'0.6032474722945804'
$arr[16,0] = 'LU_Comm_res'
$arr[16,1] = [double]'1510.326954918995'
$arr[16,2] = [double]'0.7440091290750337'
$arr[17,0] = 'LU_UrbFab_res'
$arr[17,1] = [double]'-5389.889595015151'
$arr[17,2] = [double]'0.2332818675747863'
$arr[18,0] = 'bike_lane_share_res'
$arr[18,1] = [double]'-3507.159299231108'
$arr[18,2] = [double]'0.4764981486502424'
$ws.Range("A2:C20").Value = $arr
